# "add exercise session to calender"
# Update the 30.01 lecture-plan entry in column D (row 5) to advertise the
# new Exercise session, wrap the text, size the row to fit it, and leave
# the selection where the author left it after editing (D6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New wording for the 30.01 session: "Practical session" -> "Exercise session"
# (with a line break before "Exercise session", matching the authored cell).
$ws.Cells.Item(5, 4).Value = "30.01: <strong>`nExercise session</strong> in Aud J"

# Wrap the text so the two lines are visible, and size row 5 to fit it.
$ws.Cells.Item(5, 4).WrapText = $true
$ws.Rows.Item(5).RowHeight = 30

# Leave the cursor on D6, matching the post-edit selection.
$ws.Range("D6").Select() | Out-Null
